$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 41
    3  = 40
    4  = 6
    5  = 24
    6  = 16.5
    7  = 19
    8  = 30.9090909090909
    9  = 36.3636363636364
    10 = 25.4545454545455
    11 = 42.7272727272727
    12 = 29
    13 = 25
    14 = 34.5454545454545
    15 = 24
    16 = 32.7272727272727
    17 = 22
    18 = 30.9090909090909
    19 = 48.1818181818182
    20 = 14.5454545454545
    21 = 27.2727272727273
    22 = 50
    23 = 14.5454545454545
}

foreach ($row in $values.Keys) {
    $ws.Range("AF$row").Value = $values[$row]
}

$ws.Range("AL15").Select()
